# Regenerate the "K" column (column G) values for the weems_jordan save-data
# sheet. These values are recalculated from the underlying box-score data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals") - only the K totals per game (rows 2-37) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 1
    4  = 3
    5  = 1
    6  = 0
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 4
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 3
    26 = 3
    27 = 1
    28 = 0
    29 = 1
    30 = 3
    31 = 0
    32 = 0
    33 = 2
    34 = 3
    35 = 1
    36 = 1
    37 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
